# Commit: "added PMID and notes columns to data structure"
#
# - "studies" sheet gains a new "PMID" column (H1), right after "reference_year".
# - "counts" sheet gains a new "notes" column (F1), right after "total_num".
# - The "counts" sheet becomes the active sheet / selection lands on the
#   newly-added header's first data cell in each sheet (matches how Excel
#   leaves the cursor on the cell you just typed into and which sheet you
#   were last on when the file was saved).

$wb = $excel.ActiveWorkbook

# --- studies sheet: add "PMID" header column ---------------------------
$studies = $wb.Worksheets.Item("studies")
$studies.Range("H1").Value = "PMID"
$studies.Activate()
$studies.Range("H2").Select()

# --- counts sheet: add "notes" header column ----------------------------
$counts = $wb.Worksheets.Item("counts")
$counts.Range("F1").Value = "notes"
$counts.Activate()
$counts.Range("F2").Select()
